$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.634.30'
$ws.Range("E2").Value = '  +6.81%  '

$ws.Range("D3").Value = '3.000.46'
$ws.Range("E3").Value = '  +3.81%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Formula = "'585.33"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.09%  '

$ws.Range("D6").Formula = "'153.77"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +6.75%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").Value = '2.997.03'
$ws.Range("E8").Value = '  +3.72%  '

$ws.Range("B9").Value = 'XRP'
$ws.Range("C9").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D9").Formula = "'0.516"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.08%  '

$ws.Range("D10").Formula = "'6.97"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.41%  '

$ws.Range("D11").Formula = "'0.153"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +4.33%  '

$ws.Range("E12").Value = '  +3.75%  '

$ws.Range("E13").Value = '  +3.11%  '

$ws.Range("D14").Formula = "'33.97"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +5.99%  '

$ws.Range("E15").Value = '  +0.62%  '

$ws.Range("D16").Value = '65.557.37'
$ws.Range("E16").Value = '  +6.67%  '

$ws.Range("D17").Value = '3.497.65'
$ws.Range("E17").Value = '  +3.82%  '

$ws.Range("D18").Formula = "'6.93"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +5.75%  '

$ws.Range("D19").Value = '2.998.53'
$ws.Range("E19").Value = '  +3.89%  '

$ws.Range("D20").Formula = "'452.88"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +4.59%  '

$ws.Range("D21").Formula = "'13.71"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +4.65%  '

$ws.Range("D22").Formula = "'0.681"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +3.76%  '

$ws.Range("D23").Formula = "'7.33"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +7.11%  '

$ws.Range("D24").Formula = "'81.36"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.61%  '

$ws.Range("D25").Formula = "'12.45"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +3.50%  '

$ws.Range("D26").Formula = "'2.23"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +10.60%  '

$ws.Range("D27").Formula = "'10.66"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +6.59%  '

$ws.Range("E28").Value = '  -0.04%  '

$ws.Range("E29").Value = '  +17.02%  '

$ws.Range("E30").Value = '  +10.98%  '

$ws.Range("E31").Value = '  +3.84%  '

$ws.Range("E32").Value = '  -2.05%  '

$ws.Range("D33").Formula = "'26.90"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +5.43%  '

$ws.Range("E34").Value = '  +3.84%  '

$ws.Range("D35").Formula = "'0.998"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.26%  '

$ws.Range("D36").Formula = "'0.988"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +2.91%  '

$ws.Range("D37").Formula = "'5.77"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +6.87%  '

$ws.Range("D38").Formula = "'2.12"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +9.49%  '

$ws.Range("D39").Formula = "'45.74"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +17.07%  '

$ws.Range("D40").Formula = "'49.20"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.66%  '

$ws.Range("E41").Value = '  +2.55%  '

$ws.Range("E42").Value = '  +5.97%  '

$ws.Range("D43").Formula = "'0.299"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +11.84%  '

$ws.Range("E44").Value = '  +2.22%  '

$ws.Range("D45").Formula = "'386.65"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +12.30%  '

$ws.Range("D46").Value = '2.766.48'
$ws.Range("E46").Value = '  +2.11%  '

$ws.Range("E47").Value = '  +4.78%  '

$ws.Range("E48").Value = '  +1.64%  '

$ws.Range("E49").Value = '  -0.02%  '

$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Formula = "'23.25"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +7.74%  '

$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").Formula = "'0.106"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.76%  '
